{"js": "// Update the date heading paragraph.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nif (paras.items.length > 0) {\n  const firstPara = paras.items[0];\n  firstPara.load(\"text\");\n  await context.sync();\n  if (firstPara.text.trim() === \"2024-03-09 Saturday\") {\n    firstPara.insertText(\"2024-03-10 Sunday\", Word.InsertLocation.replace);\n  }\n}\n\n// Update the arithmetic table: replace each cell's text in row-major order\n// (20 rows x 5 columns) with the new problem for that position.\nconst newValues = [\n    [\"30-8=\", \"58-26=\", \"9+86=\", \"76-37=\", \"8+47=\"],\n    [\"51+20=\", \"58+34=\", \"35-20=\", \"73-33=\", \"60+38=\"],\n    [\"25+54=\", \"26+50=\", \"88-66=\", \"34+33=\", \"23+52=\"],\n    [\"62+14=\", \"90-42=\", \"96-41=\", \"24+54=\", \"0+5=\"],\n    [\"53+4=\", \"60+0=\", \"8+54=\", \"67-34=\", \"7+66=\"],\n    [\"83-50=\", \"10+6=\", \"63+15=\", \"15-2=\", \"72-29=\"],\n    [\"94-78=\", \"16+31=\", \"1+81=\", \"22+68=\", \"66+22=\"],\n    [\"0+80=\", \"51+7=\", \"9+8=\", \"98-2=\", \"24-23=\"],\n    [\"17+8=\", \"72+0=\", \"36+48=\", \"15+15=\", \"71-1=\"],\n    [\"56+1=\", \"35+50=\", \"69-23=\", \"24+11=\", \"23+63=\"],\n    [\"9+85=\", \"40+14=\", \"45-8=\", \"32+19=\", \"47-24=\"],\n    [\"21+2=\", \"15+15=\", \"28+57=\", \"62+12=\", \"2+0=\"],\n    [\"73-62=\", \"65+27=\", \"41-38=\", \"53-39=\", \"98-73=\"],\n    [\"68-6=\", \"4+69=\", \"12+3=\", \"89-5=\", \"67+8=\"],\n    [\"89-81=\", \"63-33=\", \"47-22=\", \"93-70=\", \"81-38=\"],\n    [\"25+14=\", \"43-26=\", \"58-21=\", \"16+24=\", \"1+25=\"],\n    [\"49-12=\", \"7+79=\", \"65-60=\", \"9+70=\", \"86-46=\"],\n    [\"21+5=\", \"8+26=\", \"22+71=\", \"18+46=\", \"23+29=\"],\n    [\"92-6=\", \"51+13=\", \"52+39=\", \"63+34=\", \"24+29=\"],\n    [\"38+43=\", \"11+81=\", \"81+15=\", \"84-13=\", \"55-38=\"],\n  ];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.values = newValues;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph).\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd([char]13, [char]7) -eq \"2024-03-09 Saturday\") {\n    $titlePara.Range.Text = \"2024-03-10 Sunday\"\n}\n\n# Build the row-major (20 rows x 5 columns) replacement grid for the\n# arithmetic table, matching the target edit cell-for-cell.\n$newValues = @(\n    @(\"30-8=\", \"58-26=\", \"9+86=\", \"76-37=\", \"8+47=\"),\n    @(\"51+20=\", \"58+34=\", \"35-20=\", \"73-33=\", \"60+38=\"),\n    @(\"25+54=\", \"26+50=\", \"88-66=\", \"34+33=\", \"23+52=\"),\n    @(\"62+14=\", \"90-42=\", \"96-41=\", \"24+54=\", \"0+5=\"),\n    @(\"53+4=\", \"60+0=\", \"8+54=\", \"67-34=\", \"7+66=\"),\n    @(\"83-50=\", \"10+6=\", \"63+15=\", \"15-2=\", \"72-29=\"),\n    @(\"94-78=\", \"16+31=\", \"1+81=\", \"22+68=\", \"66+22=\"),\n    @(\"0+80=\", \"51+7=\", \"9+8=\", \"98-2=\", \"24-23=\"),\n    @(\"17+8=\", \"72+0=\", \"36+48=\", \"15+15=\", \"71-1=\"),\n    @(\"56+1=\", \"35+50=\", \"69-23=\", \"24+11=\", \"23+63=\"),\n    @(\"9+85=\", \"40+14=\", \"45-8=\", \"32+19=\", \"47-24=\"),\n    @(\"21+2=\", \"15+15=\", \"28+57=\", \"62+12=\", \"2+0=\"),\n    @(\"73-62=\", \"65+27=\", \"41-38=\", \"53-39=\", \"98-73=\"),\n    @(\"68-6=\", \"4+69=\", \"12+3=\", \"89-5=\", \"67+8=\"),\n    @(\"89-81=\", \"63-33=\", \"47-22=\", \"93-70=\", \"81-38=\"),\n    @(\"25+14=\", \"43-26=\", \"58-21=\", \"16+24=\", \"1+25=\"),\n    @(\"49-12=\", \"7+79=\", \"65-60=\", \"9+70=\", \"86-46=\"),\n    @(\"21+5=\", \"8+26=\", \"22+71=\", \"18+46=\", \"23+29=\"),\n    @(\"92-6=\", \"51+13=\", \"52+39=\", \"63+34=\", \"24+29=\"),\n    @(\"38+43=\", \"11+81=\", \"81+15=\", \"84-13=\", \"55-38=\"),\n)\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n\n$d.Save()\n"}
